$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule")

# --- Activate the Schedule sheet (this moves workbook.xml's activeTab,
#     sets tabSelected on this sheet's view, and clears it on the
#     previously active "Status" sheet) ---
$ws.Activate()

# --- Selection / freeze-pane state on the Schedule sheet ---
# Target: pane ySplit=1 topLeftCell=A2 activePane=bottomLeft state=frozen
#         selection pane=bottomLeft activeCell=I11 sqref=I11
$ws.Range("I11").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- Column widths (B:E) ---
$ws.Columns.Item(2).ColumnWidth = 36
$ws.Columns.Item(3).ColumnWidth = 10.666666666666666
$ws.Columns.Item(4).ColumnWidth = 26.166666666666668
$ws.Columns.Item(5).ColumnWidth = 27.5

# --- Row 1 height back to default (removes explicit row height) ---
$ws.Rows.Item(1).AutoFit()

# --- Schedule.Health.Standard (column I) value updates ---
$ws.Cells.Item(10, 9).Value = 4.5
$ws.Cells.Item(26, 9).Value = 4.5
$ws.Cells.Item(27, 9).Value = 4.5
$ws.Cells.Item(28, 9).Value = 4.5
$ws.Cells.Item(29, 9).Value = 4.5
$ws.Cells.Item(30, 9).Value = 4.5
$ws.Cells.Item(31, 9).Value = 4.5
$ws.Cells.Item(32, 9).Value = 4.5
$ws.Cells.Item(35, 9).Value = 4.5
$ws.Cells.Item(49, 9).Value = 4.5

# --- Row 51: Actual_date moves a year later, Schedule.Health.Standard -> 3 ---
$ws.Cells.Item(51, 7).Value = "2/19/2021"
$ws.Cells.Item(51, 9).Value = 3

# --- Schedule.Health.Standard 6 -> 3 / "completed" ---
$ws.Cells.Item(59, 9).Value = 3
$ws.Cells.Item(60, 9).Value = 3
$ws.Cells.Item(61, 9).Value = "completed"
$ws.Cells.Item(62, 9).Value = "completed"
$ws.Cells.Item(63, 9).Value = 3
$ws.Cells.Item(64, 9).Value = 3
$ws.Cells.Item(65, 9).Value = "completed"
$ws.Cells.Item(66, 9).Value = "completed"

# --- Row 67 height becomes a custom (non-default) height ---
$ws.Rows.Item(67).RowHeight = 29.25

$ws.Cells.Item(67, 9).Value = 3
$ws.Cells.Item(68, 9).Value = 3
